# Atualização de bases das ligas, do dia: 20-02-2024 às 23:00
#
# The underlying data source re-ordered two pairs of already-recorded
# fixtures (rows 78/79 and rows 82/83 swapped places) and refreshed the
# odds for three upcoming fixtures (rows 184, 185, 187) - row 184's match
# has since been played, so its result (FTHG/FTAG/FTR) and closing
# Asian-handicap odds (PLH..PL_AhUnder) are now populated too.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rows 78 and 79: the two fixtures were swapped in the source feed ---
$row78 = $ws.Range("B78:AC78")
$row79 = $ws.Range("B79:AC79")
$v78 = $row78.Value2
$v79 = $row79.Value2
$row78.Value = $v79
$row79.Value = $v78

# --- Rows 82 and 83: same kind of swap ---
$row82 = $ws.Range("B82:AC82")
$row83 = $ws.Range("B83:AC83")
$v82 = $row82.Value2
$v83 = $row83.Value2
$row82.Value = $v83
$row83.Value = $v82

# --- Row 184: match has been played, fill in the result + refresh odds ---
$ws.Range("H184").Value = 1
$ws.Range("I184").Value = 2
$ws.Range("J184").Value = "A"
$ws.Range("N184").Value = 1.75
$ws.Range("P184").Value = 4
$ws.Range("W184").Value = -1
$ws.Range("X184").Value = -1
$ws.Range("Y184").Value = 3
$ws.Range("Z184").Value = -1
$ws.Range("AA184").Value = 0.9750000000000001
$ws.Range("AB184").Value = 0.8999999999999999
$ws.Range("AC184").Value = -1

# --- Row 185: odds refresh only (match not played yet) ---
$ws.Range("R185").Value = 1.85
$ws.Range("S185").Value = 1.95
$ws.Range("T185").Value = 3
$ws.Range("U185").Value = 1.975
$ws.Range("V185").Value = 1.825

# --- Row 187: odds refresh only (match not played yet) ---
$ws.Range("N187").Value = 1.285
$ws.Range("O187").Value = 5.5
$ws.Range("P187").Value = 8
$ws.Range("R187").Value = 2
$ws.Range("S187").Value = 1.8
$ws.Range("U187").Value = 2
$ws.Range("V187").Value = 1.8
